$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.439.95"
$ws.Range("E2").Value = "  -0.82%  "
$ws.Range("D3").Value = "1.892.11"
$ws.Range("E3").Value = "  +0.09%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4900"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.24%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2937"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.27%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06700"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.20%  "
$ws.Range("D10").Value = "1.874.99"
$ws.Range("E10").Value = "  -0.78%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "16.98"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.35%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07351"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.53%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.132"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.34%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "87.87"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.78%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6637"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.36%  "
$ws.Range("D16").Value = "30.417.53"
$ws.Range("E16").Value = "  -0.73%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.44"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.96%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007819"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.000"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.01%  "
$ws.Range("D20").Value = "2.139.78"
$ws.Range("E20").Value = "  +0.21%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.306"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +11.68%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.003"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "190.21"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.136"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.95%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.474"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.67%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.64"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.34%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.24"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.48%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.930"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.25%  "
$ws.Range("E29").Value = "  +4.46%  "
$ws.Range("E30").Value = "  +2.31%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09150"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.045"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.47%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05208"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.08%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7404"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.77%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.099"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.717"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01813"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.66%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.674"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.18%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9203"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.52%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.033"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.94%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4399"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.63%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.939"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.26%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "106.18"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.36%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9924"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.74%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "68.32"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +19.25%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1372"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.23%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.572"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.75%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.018"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.89%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "34.97"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.02%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05821"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3935"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.97%  "
